$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (before the current HYD row), shifting
# existing rows 7-14 down to 8-15, to keep the BrowseProduct list sorted
# alphabetically with the new "HY3" band-area parameters entry.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the HY3 band-area parameters.
$ws.Range("A7").Value = "HY3"
$ws.Range("B7").Value = "BA1200"
$ws.Range("C7").Value = "BA1450"
$ws.Range("D7").Value = "BA1900"

# The hidden _FilterDatabase defined name tracked the old used range
# (A1:D14); grow it to match the new used range (A1:D15).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$D`$15"
    }
}

# Match the saved selection state.
$ws.Range("A8").Select()
